{"js": "// Exercise 2 explanation paragraph was rewritten based on new training logs\n// (commit: \"Uitleg ex2 aangepast o.b.v. nieuwe logs, feedback toegevoegd\").\n// Find the paragraph that starts with the old \"Looking at the accuracy plot\"\n// text (it also contained a REF field to \"Figure 1\" and was highlighted\n// yellow) and replace its whole content with the new, un-highlighted text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the target paragraph robustly instead of relying on a fixed index.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Looking at the accuracy plot\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Looking at the accuracy plot' paragraph\");\n}\n\n// New wording (replaces the old field-based \"(Figure 1)...\" sentence and the\n// yellow highlight that marked it as a draft / to-be-reviewed note).\nconst newText =\n  \"Looking at the accuracy plot (Figure 2), it seems to be the case that the \" +\n  \"performance of model 3 does not increase with training. It\\u2019s final value \" +\n  \"(0.5163) is lower than it\\u2019s starting value (0.535).  For model 1, the \" +\n  \"validation accuracy seems to have stopped increasing after 10 epochs.  The \" +\n  \"validation accuracy is not yet decreasing significantly, so overtraining has \" +\n  \"not yet occurred. The validation and training accuracy of model 3 are still \" +\n  \"increasing after 3 epochs. More training seems to be needed to improve the \" +\n  \"model performance.\";\n\n// Wipe out the existing runs (plain text + REF field + highlighted runs).\ntarget.clear();\n\n// Insert the new text as the paragraph's content.\nconst newRange = target.insertText(newText, Word.InsertLocation.start);\n\n// Make sure no highlight carries over and the language stays English (US),\n// matching the rest of the document's runs.\nnewRange.font.highlightColor = null;\nnewRange.languageId = \"en-US\";\n\nawait context.sync();\n", "ps1": "# Exercise 2 explanation paragraph was rewritten based on new training logs\n# (commit: \"Uitleg ex2 aangepast o.b.v. nieuwe logs, feedback toegevoegd\").\n# Find the paragraph that starts with the old \"Looking at the accuracy plot\"\n# text (it also contained a REF field pointing at \"Figure 1\" and was\n# highlighted yellow) and replace its whole content with the new,\n# un-highlighted text.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Looking at the accuracy plot*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Looking at the accuracy plot' paragraph\"\n}\n\n$rightQuote = [char]0x2019\n\n$newText = \"Looking at the accuracy plot (Figure 2), it seems to be the case that \" + `\n    \"the performance of model 3 does not increase with training. It\" + $rightQuote + `\n    \"s final value (0.5163) is lower than it\" + $rightQuote + \"s starting value \" + `\n    \"(0.535).  For model 1, the validation accuracy seems to have stopped \" + `\n    \"increasing after 10 epochs.  The validation accuracy is not yet decreasing \" + `\n    \"significantly, so overtraining has not yet occurred. The validation and \" + `\n    \"training accuracy of model 3 are still increasing after 3 epochs. More \" + `\n    \"training seems to be needed to improve the model performance.\"\n\n# Remove the existing runs (plain text + REF field + highlighted runs) but\n# keep the paragraph mark itself intact, then insert the new wording.\n$paraRange = $target.Range\n$textRange = $d.Range($paraRange.Start, $paraRange.End - 1)\n$textRange.Delete()\n$target.Range.InsertAfter($newText)\n\n# Make sure no highlight carries over and the language stays English (US),\n# matching the rest of the document's runs.\n$target.Range.HighlightColorIndex = 0\n$target.Range.LanguageID = \"en-US\"\n"}
